{"js": "// Replace the date line and each two-digit multiplication problem with its\n// new value. Every original text run in this document is unique, so an\n// exact, case-sensitive search-and-replace on each pair is unambiguous.\nconst replacements = [\n  [\"2025-06-04 Wednesday\", \"2025-06-05 Thursday\"],\n  [\"85\u00d714=\", \"80\u00d752=\"],\n  [\"56\u00d787=\", \"86\u00d744=\"],\n  [\"68\u00d736=\", \"22\u00d721=\"],\n  [\"31\u00d741=\", \"91\u00d788=\"],\n  [\"60\u00d736=\", \"48\u00d715=\"],\n  [\"16\u00d716=\", \"26\u00d719=\"],\n  [\"89\u00d766=\", \"13\u00d734=\"],\n  [\"43\u00d743=\", \"93\u00d721=\"],\n  [\"18\u00d715=\", \"53\u00d799=\"],\n  [\"79\u00d722=\", \"33\u00d711=\"],\n  [\"82\u00d751=\", \"50\u00d721=\"],\n  [\"92\u00d777=\", \"19\u00d749=\"],\n  [\"63\u00d773=\", \"20\u00d749=\"],\n  [\"29\u00d787=\", \"17\u00d786=\"],\n  [\"38\u00d713=\", \"82\u00d711=\"],\n  [\"17\u00d757=\", \"51\u00d782=\"],\n  [\"22\u00d730=\", \"17\u00d781=\"],\n  [\"56\u00d788=\", \"80\u00d732=\"],\n  [\"68\u00d770=\", \"42\u00d773=\"],\n  [\"24\u00d728=\", \"86\u00d782=\"],\n  [\"62\u00d781=\", \"62\u00d714=\"],\n  [\"99\u00d740=\", \"91\u00d732=\"],\n  [\"93\u00d784=\", \"24\u00d773=\"],\n  [\"66\u00d780=\", \"42\u00d725=\"],\n  [\"27\u00d757=\", \"86\u00d772=\"],\n];\n\nconst searchOptions = { matchCase: true, matchWholeWord: false };\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication problem with its\n# new value. Every original text run in this document is unique, so a\n# single Find/Replace pass per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-06-04 Wednesday\", \"2025-06-05 Thursday\"),\n  @(\"85\u00d714=\", \"80\u00d752=\"),\n  @(\"56\u00d787=\", \"86\u00d744=\"),\n  @(\"68\u00d736=\", \"22\u00d721=\"),\n  @(\"31\u00d741=\", \"91\u00d788=\"),\n  @(\"60\u00d736=\", \"48\u00d715=\"),\n  @(\"16\u00d716=\", \"26\u00d719=\"),\n  @(\"89\u00d766=\", \"13\u00d734=\"),\n  @(\"43\u00d743=\", \"93\u00d721=\"),\n  @(\"18\u00d715=\", \"53\u00d799=\"),\n  @(\"79\u00d722=\", \"33\u00d711=\"),\n  @(\"82\u00d751=\", \"50\u00d721=\"),\n  @(\"92\u00d777=\", \"19\u00d749=\"),\n  @(\"63\u00d773=\", \"20\u00d749=\"),\n  @(\"29\u00d787=\", \"17\u00d786=\"),\n  @(\"38\u00d713=\", \"82\u00d711=\"),\n  @(\"17\u00d757=\", \"51\u00d782=\"),\n  @(\"22\u00d730=\", \"17\u00d781=\"),\n  @(\"56\u00d788=\", \"80\u00d732=\"),\n  @(\"68\u00d770=\", \"42\u00d773=\"),\n  @(\"24\u00d728=\", \"86\u00d782=\"),\n  @(\"62\u00d781=\", \"62\u00d714=\"),\n  @(\"99\u00d740=\", \"91\u00d732=\"),\n  @(\"93\u00d784=\", \"24\u00d773=\"),\n  @(\"66\u00d780=\", \"42\u00d725=\"),\n  @(\"27\u00d757=\", \"86\u00d772=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Text = $oldText\n  $range.Find.Replacement.Text = $newText\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
